# Switch license from BY-NC to BY-SA on the "license / credits" slide.
#
# This slide is the 2nd slide in the deck (p:sldId id="267", r:id="rId3",
# backed by ppt/slides/slide2.xml).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- 1. Text edits -------------------------------------------------------
# "CC BY-NC 4.0. To view a copy of this license, visit " -> "CC BY-SA 4.0. ..."
$full = $tr.Text
$idx = $full.IndexOf("BY-NC ")
$run = $tr.Characters($idx + 1, 6)
$run.Text = "BY-SA "

# "https://creativecommons.org/licenses/by-nc/4.0" -> ".../by-sa/4.0"
$full = $tr.Text
$idx = $full.IndexOf("creativecommons.org/licenses/by-nc/4.0")
$run = $tr.Characters($idx + 1, 39)
$run.Text = "creativecommons.org/licenses/by-sa/4.0"

# --- 2. Nudge the title placeholder's position ----------------------------
$shp.Left = 566057 / 12700.0

# --- 3. Drop the vestigial empty <p:timing> block -------------------------
# Adding then immediately deleting an animation effect clears out the
# leftover empty timeline node entirely.
$seq = $s.TimeLine.MainSequence
$effect = $seq.AddEffect($shp, 1, 0, 1)
$effect.Delete()
